# Commit: "updated uart1 wire color- yellow and white are swapped at
# other connector end"
#
# On the (single) slide, two small wire-color label textboxes are
# expanded from a bare color abbreviation into a "<abbr> or <abbr>?"
# question, since at the far end of the UART1 cable the yellow/white
# wires are swapped:
#   id=75 "yel" -> "yel or wht?"
#   id=77 "wht" -> "wht or yel?"
# Both textboxes also grow (wrap="none" + spAutoFit, so the shape is
# resized/repositioned to fit the new, longer caption).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$yelShape = $null
$whtShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 75) { $yelShape = $shp }
    if ($shp.Id -eq 77) { $whtShape = $shp }
}

# ---- "yel" textbox -> "yel or wht?" ----------------------------------
$tr = $yelShape.TextFrame.TextRange

# Split the original "yel" run into "y" | "el" (re-assigning the first
# character forces PowerPoint to split the run at that boundary).
$tr.Characters(1, 1).Text = "y"

# Append the new tail text as its own run.
[void]$tr.InsertAfter(" or wht?")

# Split " or wht?" into " or " | "wht" | "?" by re-assigning the "wht"
# sub-range (re-assigning a middle slice splits runs on both sides).
$tr.Characters(8, 3).Text = "wht"

# Resize/reposition the textbox (autosize grew it to fit "yel or wht?").
$yelShape.Left = 516.81132
$yelShape.Top = 441.0771
$yelShape.Width = 96.9527
$yelShape.Height = 29.0813

# ---- "wht" textbox -> "wht or yel?" ----------------------------------
$tr2 = $whtShape.TextFrame.TextRange

# Split the original "wht" run into "w" | "ht".
$tr2.Characters(1, 1).Text = "w"

# Append the new tail text as its own run.
[void]$tr2.InsertAfter(" or yel?")

# Split " or yel?" into " or " | "yel" | "?".
$tr2.Characters(8, 3).Text = "yel"

# Resize/reposition the textbox (autosize grew it to fit "wht or yel?").
$whtShape.Left = 516.81144
$whtShape.Top = 416.69222
$whtShape.Width = 96.9527
$whtShape.Height = 29.0813
